$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80 (hunk 0)
$ws.Range("H80").Value = 27778536
$ws.Range("I80").Value = 466.66666
$ws.Range("J80").Value = 41667572
$ws.Range("K80").Value = 1399.99998
$ws.Range("L80").Value = 125002716
$ws.Range("M80").Value = -401.9999800000001
$ws.Range("N80").Value = -125004712

# Row 83 (hunk 1)
$ws.Range("H83").Value = 27778536
$ws.Range("I83").Value = 466.66666
$ws.Range("J83").Value = 41667572
$ws.Range("K83").Value = 4199.99994
$ws.Range("L83").Value = 375008148
$ws.Range("M83").Value = 792.0000600000003
$ws.Range("N83").Value = -375018132

# Row 88 (hunk 2)
$ws.Range("H88").Value = 714
$ws.Range("I88").Value = 954.5
$ws.Range("K88").Value = 954.5
$ws.Range("M88").Value = -548.5

# Row 91 (hunk 3)
$ws.Range("H91").Value = 714
$ws.Range("I91").Value = 954.5
$ws.Range("K91").Value = 954.5
$ws.Range("M91").Value = 449.5

# Row 94 (hunk 4)
$ws.Range("H94").Value = 1698
$ws.Range("I94").Value = 1698
$ws.Range("K94").Value = 1698
$ws.Range("M94").Value = -1247

# Row 99 (hunk 5)
$ws.Range("H99").Value = 90909460
$ws.Range("I99").Value = 456.125
$ws.Range("K99").Value = 1368.375
$ws.Range("M99").Value = 129.625

# Row 103 (hunk 6)
$ws.Range("H103").Value = 437.66666
$ws.Range("I103").Value = 442.5
$ws.Range("J103").Value = 399
$ws.Range("K103").Value = 1327.5
$ws.Range("L103").Value = 1197
$ws.Range("M103").Value = -741.5
$ws.Range("N103").Value = -2369

# Row 106 (hunk 7)
$ws.Range("H106").Value = 83335630
$ws.Range("I106").Value = 86958730
$ws.Range("K106").Value = 86958730
$ws.Range("M106").Value = -86958099

# Row 137 (hunk 8)
$ws.Range("H137").Value = 5792.9346
$ws.Range("I137").Value = 4312.7417
$ws.Range("K137").Value = 12938.2251
$ws.Range("M137").Value = -10388.2251

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (hunk 9)
$ws.Range("H2").Value = 2049.5
$ws.Range("I2").Value = 2000
$ws.Range("J2").Value = 2099
$ws.Range("K2").Value = 2000
$ws.Range("L2").Value = 2099
$ws.Range("M2").Value = -1887
$ws.Range("N2").Value = -2325

# Row 32 (hunk 10)
$ws.Range("H32").Value = 177966.1
$ws.Range("I32").Value = 201283.16
$ws.Range("J32").Value = 11415.571
$ws.Range("K32").Value = 201283.16
$ws.Range("L32").Value = 11415.571
$ws.Range("M32").Value = -200996.16
$ws.Range("N32").Value = -11989.571

# Row 74 (hunk 11)
$ws.Range("H74").Value = 24945.75
$ws.Range("I74").Value = 2171.6177
$ws.Range("K74").Value = 2171.6177
$ws.Range("M74").Value = -1297.6177

# Row 77 (hunk 12)
$ws.Range("H77").Value = 24945.75
$ws.Range("I77").Value = 2171.6177
$ws.Range("K77").Value = 10858.0885
$ws.Range("M77").Value = -6490.088499999998

# Row 102 (hunk 13)
$ws.Range("H102").Value = 2333.3333
$ws.Range("I102").Value = 2000
$ws.Range("K102").Value = 2000
$ws.Range("M102").Value = -378

# Row 110 (hunk 14)
$ws.Range("H110").Value = 1307.5385
$ws.Range("I110").Value = 1166.5834
$ws.Range("J110").Value = 2999
$ws.Range("K110").Value = 1166.5834
$ws.Range("L110").Value = 2999
$ws.Range("M110").Value = 878.4166
$ws.Range("N110").Value = -7089

# Row 116 (hunk 15)
$ws.Range("H116").Value = 2049.5
$ws.Range("I116").Value = 2000
$ws.Range("J116").Value = 2099
$ws.Range("K116").Value = 2000
$ws.Range("L116").Value = 2099
$ws.Range("M116").Value = 294
$ws.Range("N116").Value = -6687

# Row 122 (hunk 16)
$ws.Range("H122").Value = 1944.4
$ws.Range("I122").Value = 993
$ws.Range("K122").Value = 2979
$ws.Range("M122").Value = -529

# Row 132 (hunk 17)
$ws.Range("H132").Value = 1919.279
$ws.Range("I132").Value = 1330.3667
$ws.Range("K132").Value = 3991.1001
$ws.Range("M132").Value = -1461.1001

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (hunk 18)
$ws.Range("H3").Value = 2049.5
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 2099
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 2099
$ws.Range("M3").Value = -1886
$ws.Range("N3").Value = -2327

# Row 99 (hunk 19)
$ws.Range("H99").Value = 9888.706
$ws.Range("I99").Value = 12316
$ws.Range("K99").Value = 12316
$ws.Range("M99").Value = -10818

# Row 105 (hunk 20)
$ws.Range("H105").Value = 9482.147999999999
$ws.Range("I105").Value = 9123.25
$ws.Range("K105").Value = 9123.25
$ws.Range("M105").Value = -7376.25

# Row 113 (hunk 21)
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

# Row 134 (hunk 22)
$ws.Range("H134").Value = 1495.5306
$ws.Range("I134").Value = 1341.9412
$ws.Range("K134").Value = 4025.8236
$ws.Range("M134").Value = -1490.8236

$ws = $wb.Worksheets.Item("CRP")
# Row 99 (hunk 23)
$ws.Range("H99").Value = 10000400
$ws.Range("I99").Value = 10000400
$ws.Range("K99").Value = 10000400
$ws.Range("M99").Value = -9998902

# Row 105 (hunk 24)
$ws.Range("H105").Value = 4799
$ws.Range("I105").Value = 4665
$ws.Range("K105").Value = 4665
$ws.Range("M105").Value = -2918

# Row 126 (hunk 25)
$ws.Range("H126").Value = 10000400
$ws.Range("I126").Value = 10000400
$ws.Range("K126").Value = 30001200
$ws.Range("M126").Value = -29998730

# Row 132 (hunk 26)
$ws.Range("H132").Value = 1917.6666
$ws.Range("J132").Value = 2486.7778
$ws.Range("L132").Value = 7460.3334
$ws.Range("N132").Value = -12520.3334

# Row 141 (hunk 27)
$ws.Range("H141").Value = 286378
$ws.Range("I141").Value = 77321
$ws.Range("K141").Value = 77321
$ws.Range("M141").Value = -72141

$ws = $wb.Worksheets.Item("CUL")
# Row 97 (hunk 28)
$ws.Range("H97").Value = 616
$ws.Range("I97").Value = 498.33334
$ws.Range("J97").Value = 792.5
$ws.Range("K97").Value = 1495.00002
$ws.Range("L97").Value = 2377.5
$ws.Range("M97").Value = -999.0000199999999
$ws.Range("N97").Value = -3369.5

# Row 122 (hunk 29)
$ws.Range("H122").Value = 3704700.5
$ws.Range("I122").Value = 4167399.5
$ws.Range("K122").Value = 37506595.5
$ws.Range("M122").Value = -37504145.5

# Row 138 (hunk 30)
$ws.Range("H138").Value = 3069.3684
$ws.Range("I138").Value = 3153.611
$ws.Range("K138").Value = 9460.832999999999
$ws.Range("M138").Value = -4320.832999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (hunk 31)
$ws.Range("H70").Value = 6338.6665
$ws.Range("I70").Value = 6008
$ws.Range("J70").Value = 7000
$ws.Range("K70").Value = 6008
$ws.Range("L70").Value = 7000
$ws.Range("M70").Value = -5738
$ws.Range("N70").Value = -7540

# Row 73 (hunk 32)
$ws.Range("H73").Value = 6338.6665
$ws.Range("I73").Value = 6008
$ws.Range("J73").Value = 7000
$ws.Range("K73").Value = 6008
$ws.Range("L73").Value = 7000
$ws.Range("M73").Value = -5072
$ws.Range("N73").Value = -8872

# Row 97 (hunk 33)
$ws.Range("H97").Value = 100596.6
$ws.Range("I97").Value = 71827.28999999999
$ws.Range("K97").Value = 71827.28999999999
$ws.Range("M97").Value = -71331.28999999999

# Row 102 (hunk 34)
$ws.Range("H102").Value = 17242830
$ws.Range("I102").Value = 20834688
$ws.Range("K102").Value = 20834688
$ws.Range("M102").Value = -20833066

# Row 113 (hunk 35)
$ws.Range("H113").Value = 6000
$ws.Range("I113").Value = 6000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 6000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -3830
$ws.Range("N113").ClearContents()

# Row 132 (hunk 36)
$ws.Range("H132").Value = 14031.777
$ws.Range("I132").Value = 23091.6
$ws.Range("J132").Value = 2707
$ws.Range("K132").Value = 69274.79999999999
$ws.Range("L132").Value = 8121
$ws.Range("M132").Value = -66744.79999999999
$ws.Range("N132").Value = -13181

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (hunk 37)
$ws.Range("H16").Value = 1350
$ws.Range("I16").Value = 700
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 700
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -530
$ws.Range("N16").Value = -2340

# Row 40 (hunk 38)
$ws.Range("H40").Value = 2764.0908
$ws.Range("I40").Value = 3071.4285
$ws.Range("J40").Value = 2226.25
$ws.Range("K40").Value = 3071.4285
$ws.Range("L40").Value = 2226.25
$ws.Range("M40").Value = -2935.4285
$ws.Range("N40").Value = -2498.25

# Row 82 (hunk 39)
$ws.Range("H82").Value = 556.1429000000001
$ws.Range("I82").Value = 665.8570999999999
$ws.Range("K82").Value = 665.8570999999999
$ws.Range("M82").Value = -304.8570999999999

# Row 85 (hunk 40)
$ws.Range("H85").Value = 556.1429000000001
$ws.Range("I85").Value = 665.8570999999999
$ws.Range("K85").Value = 665.8570999999999
$ws.Range("M85").Value = 582.1429000000001

# Row 100 (hunk 41)
$ws.Range("H100").Value = 2885.182
$ws.Range("I100").Value = 3161.9
$ws.Range("J100").Value = 118
$ws.Range("K100").Value = 3161.9
$ws.Range("L100").Value = 118
$ws.Range("M100").Value = -2620.9
$ws.Range("N100").Value = -1200

# Row 132 (hunk 42)
$ws.Range("H132").Value = 4469.154
$ws.Range("I132").Value = 3348.7778
$ws.Range("K132").Value = 10046.3334
$ws.Range("M132").Value = -7516.3334

# Row 136 (hunk 43)
$ws.Range("H136").Value = 35489.035
$ws.Range("I136").Value = 54474.633
$ws.Range("K136").Value = 163423.899
$ws.Range("M136").Value = -160873.899

$ws = $wb.Worksheets.Item("WVR")
# Row 3 (hunk 44)
$ws.Range("H3").Value = 341699.66
$ws.Range("I3").Value = 100
$ws.Range("K3").Value = 100
$ws.Range("M3").Value = 14

# Row 92 (hunk 45)
$ws.Range("H92").Value = 75500
$ws.Range("J92").Value = 75500
$ws.Range("L92").Value = 75500
$ws.Range("N92").Value = -80492

# Row 100 (hunk 46)
$ws.Range("H100").Value = 737
$ws.Range("I100").Value = 630.44446
$ws.Range("K100").Value = 1260.88892
$ws.Range("M100").Value = -719.8889200000001
